$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.285.65"
$ws.Range("E2").Value = "  -0.76%  "
$ws.Range("D3").Value = "1.870.23"
$ws.Range("E3").Value = "  -0.41%  "
$ws.Range("D4").Value = "'0.9999"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'0.7112"
$ws.Range("E5").Value = "  -1.01%  "
$ws.Range("D6").Value = "'241.86"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "'0.3110"
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("D9").Value = "'0.07715"
$ws.Range("E9").Value = "  -2.34%  "
$ws.Range("D10").Value = "'24.72"
$ws.Range("E10").Value = "  -2.76%  "
$ws.Range("D11").Value = "'0.08402"
$ws.Range("E11").Value = "  +1.54%  "
$ws.Range("D12").Value = "1.884.30"
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").Value = "'5.224"
$ws.Range("E13").Value = "  -1.08%  "
$ws.Range("D14").Value = "'0.7122"
$ws.Range("E14").Value = "  -2.33%  "
$ws.Range("D15").Value = "'91.08"
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("D16").Value = "29.293.26"
$ws.Range("E16").Value = "  -0.76%  "
$ws.Range("D17").Value = "'0.000008155"
$ws.Range("E17").Value = "  +3.72%  "
$ws.Range("D18").Value = "'5.936"
$ws.Range("E18").Value = "  +0.56%  "
$ws.Range("D19").Value = "'243.53"
$ws.Range("E19").Value = "  -0.98%  "
$ws.Range("D20").Value = "2.124.33"
$ws.Range("E20").Value = "  -0.79%  "
$ws.Range("D21").Value = "'13.15"
$ws.Range("E21").Value = "  -1.39%  "
$ws.Range("D22").Value = "'0.9995"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").Value = "'7.864"
$ws.Range("E23").Value = "  -2.67%  "
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").Value = "'0.1628"
$ws.Range("E25").Value = "  +0.26%  "
$ws.Range("D26").Value = "'164.36"
$ws.Range("E26").Value = "  +0.58%  "
$ws.Range("D27").Value = "'9.015"
$ws.Range("E27").Value = "  -0.44%  "
$ws.Range("D28").Value = "'18.49"
$ws.Range("E28").Value = "  +0.90%  "
$ws.Range("E29").Value = "  +0.89%  "
$ws.Range("D30").Value = "'4.405"
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("E31").Value = "  -3.40%  "
$ws.Range("D32").Value = "'4.279"
$ws.Range("E32").Value = "  +4.11%  "
$ws.Range("D33").Value = "'0.05177"
$ws.Range("E33").Value = "  -0.84%  "
$ws.Range("D34").Value = "'0.7757"
$ws.Range("E34").Value = "  +6.49%  "
$ws.Range("D35").Value = "'1.915"
$ws.Range("E35").Value = "  -1.88%  "
$ws.Range("D36").Value = "'1.171"
$ws.Range("E36").Value = "  -2.31%  "
$ws.Range("D37").Value = "'2.679"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").Value = "'0.01860"
$ws.Range("E38").Value = "  -0.72%  "
$ws.Range("D39").Value = "'2.711"
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("D40").Value = "1.160.38"
$ws.Range("E40").Value = "  -3.63%  "
$ws.Range("D41").Value = "'6.397"
$ws.Range("E41").Value = "  +3.44%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'73.31"
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'0.8918"
$ws.Range("E43").Value = "  -2.09%  "
$ws.Range("D44").Value = "'0.9998"
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("D45").Value = "'103.14"
$ws.Range("E45").Value = "  +0.68%  "
$ws.Range("D46").Value = "2.020.30"
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("D47").Value = "'0.5191"
$ws.Range("E47").Value = "  -2.00%  "
$ws.Range("D48").Value = "'1.793"
$ws.Range("E48").Value = "  -0.52%  "
$ws.Range("D49").Value = "'9.392"
$ws.Range("E49").Value = "  +0.94%  "
$ws.Range("D50").Value = "'0.4299"
$ws.Range("E50").Value = "  -0.66%  "
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").Value = "'7.045"
$ws.Range("E51").Value = "  -0.51%  "

$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
